$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Age value for Carlos (row 8, column D) from 40 to 42
$ws.Range("D8").Value = 42

# Move the active selection to D8 (previously F8)
$ws.Range("D8").Select()
